$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly price-list feed shifted: the existing rows 29-56 move down one
# row (to 30-57) and a brand-new latest-week record is inserted at row 29.
$ws.Rows.Item(29).Insert()

$ws.Range("A29").Value = 1
$ws.Range("B29").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C29").Value = "Arica y Parinacota"
$ws.Range("D29").Value = 44944
$ws.Range("E29").Value = 15
$ws.Range("F29").Value = 100112028
$ws.Range("G29").Value = "Sandia"
$ws.Range("H29").Value = "Sin especificar"
$ws.Range("I29").Value = "Primera"
$ws.Range("J29").Value = 800
$ws.Range("K29").Value = 630
$ws.Range("L29").Value = 650
$ws.Range("M29").Value = 640
$ws.Range("N29").Value = "$/kilo (volumen en unidades)"
$ws.Range("O29").Value = "Perú"
$ws.Range("P29").Value = 640
$ws.Range("Q29").Value = 1
$ws.Range("R29").Value = "Hortaliza"
